# Generate Report for Handoff
#
# A new handoff xliff was generated for file
# "1700e5e4-a219-4d99-97bf-4693a0b01e11" (row 5 of the per-language detail
# sheets). This updates the "Latest Handoff Datetime" cell (column H,
# row 5) on both the "zh-cn" and "de-de" worksheets. The overview sheet's
# "Latest HO Xliff Generate Date" column (G, row 5) refers to the same
# underlying shared-string slot as the de-de sheet's value, so it is
# expected to update automatically once that value changes.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("H5").Value = "2016-08-30 00:43:10"
$wsDeDe.Range("H5").Value = "2016-08-30 00:43:15"
